$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = ' iAU_TC_ID_181'
$ws.Cells.Item(2, 2).Value = '@RegressionA Pre-Request Verify Elumina Login and Create Exam'
$ws.Cells.Item(2, 3).Value = 'passed'

$ws.Cells.Item(3, 1).Value = ' iAU_TC_ID_181.,iAU_TC_ID_183.,iAU_TC_ID_199'
$ws.Cells.Item(3, 2).Value = '@RegressionA Pre-Request "Validation of Delivery --> Add New Users"'
$ws.Cells.Item(3, 3).Value = 'passed'

$ws.Cells.Item(4, 1).Value = 'iAU_TC_ID_181'
$ws.Cells.Item(4, 2).Value = '@RegressionA Validation of Manage Delivery --> Delete Users '
$ws.Cells.Item(4, 3).Value = 'passed'

$ws.Cells.Item(5, 1).Value = ' iAU_TC_ID_190'
$ws.Cells.Item(5, 2).Value = '@RegressionA Validation of Manage Delivery--> Assign Venue and Booking Details'
$ws.Cells.Item(5, 3).Value = 'passed'

$ws.Cells.Item(6, 1).Value = 'iAU_TC_ID_182.,iAU_TC_ID_184'
$ws.Cells.Item(6, 2).Value = '@RegressionA Validation of Manage Delivery--> Delete Users (Negative Scenario) '
$ws.Cells.Item(6, 3).Value = 'passed'

$ws.Cells.Item(7, 1).Value = ' iAU_TC_ID_185.,iAU_TC_ID_186'
$ws.Cells.Item(7, 2).Value = '@RegressionA Validation of Manage Delivery --> Download User Details'
$ws.Cells.Item(7, 3).Value = 'passed'

$ws.Cells.Item(8, 1).Value = ' iAU_TC_ID_189'
$ws.Cells.Item(8, 2).Value = '@RegressionA Validation of Manage Delivery--> Generate Temp ID'
$ws.Cells.Item(8, 3).Value = 'passed'

$ws.Cells.Item(9, 1).Value = ' iAU_TC_ID_196'
$ws.Cells.Item(9, 2).Value = '@RegressionA Validation of Manage Delivery--> Reset Password'
$ws.Cells.Item(9, 3).Value = 'passed'

$ws.Cells.Item(10, 1).Value = ' iAU_TC_ID_198'
$ws.Cells.Item(10, 2).Value = '@RegressionA Validation of Manage Delivery --  Exam Administrator Manages Special Consideration'
$ws.Cells.Item(10, 3).Value = 'passed'

$ws.Cells.Item(11, 1).Value = ' iAU_TC_ID_187.,iAU_TC_ID_195.,iAU_TC_ID_200'
$ws.Cells.Item(11, 2).Value = '@RegressionA Validation of Manage Delivery --> Bulk Download User Details , Validation of Delivery --> Add New Users(Negative scenario) '
$ws.Cells.Item(11, 3).Value = 'passed'

$ws.Cells.Item(12, 1).Value = ' iAU_TC_ID_201'
$ws.Cells.Item(12, 2).Value = '@RegressionA Validation of Delivery --> Add Existing Users'
$ws.Cells.Item(12, 3).Value = 'passed'

$ws.Cells.Item(13, 1).Value = ' iAU_TC_ID_197'
$ws.Cells.Item(13, 2).Value = '@RegressionA Validation of Manage Delivery --> Bulk Candidate Response Download'
$ws.Cells.Item(13, 3).Value = 'passed'
